$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9226706027984619
$ws.Range("B1").Value = 1.661399006843567
$ws.Range("D1").Value = 1.596467018127441
$ws.Range("E1").Value = 1.043352842330933
